$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.23
$ws.Range("Q2").Value = 1.23
$ws.Range("S2").Value = 1.23
$ws.Range("F3").Value = 1.79
$ws.Range("G3").Value = 1.81
$ws.Range("I3").Value = 5.1
$ws.Range("K3").Value = 4.3
$ws.Range("Q3").Value = 1.71
$ws.Range("S3").Value = 2.8
$ws.Range("V3").Value = 1.24
$ws.Range("F4").Value = 4.8
$ws.Range("G4").Value = 6.2
$ws.Range("H4").Value = 1.87
$ws.Range("I4").Value = 2.08
$ws.Range("J4").Value = 3
$ws.Range("L4").Value = 1.56
$ws.Range("N4").Value = 2.48
$ws.Range("Q4").Value = 2.58
$ws.Range("S4").Value = 4.8
$ws.Range("U4").Value = 1.66
$ws.Range("V4").Value = 1.92
$ws.Range("W4").Value = 1.2
$ws.Range("AA4").Value = 26
$ws.Range("AC4").Value = 8.199999999999999
$ws.Range("AF4").Value = 48
$ws.Range("AG4").Value = 28
$ws.Range("AJ4").Value = 200
$ws.Range("AK4").Value = 130
$ws.Range("AL4").Value = 150
$ws.Range("AN4").Value = 230
$ws.Range("J5").Value = 3.95
$ws.Range("F6").Value = 4.3
$ws.Range("G6").Value = 5.8
$ws.Range("H6").Value = 1.8
$ws.Range("I6").Value = 2.02
$ws.Range("P6").Value = 1.88
$ws.Range("V6").Value = 1.98
$ws.Range("W6").Value = 1.21
$ws.Range("AA6").Value = 980
$ws.Range("AE6").Value = 980
$ws.Range("AF6").Value = 980
$ws.Range("AG6").Value = 980
$ws.Range("AH6").Value = 980
$ws.Range("AI6").Value = 980
$ws.Range("F7").Value = 1.27
$ws.Range("G7").Value = 1.34
$ws.Range("H7").Value = 10.5
$ws.Range("I7").Value = 15.5
$ws.Range("K7").Value = 7.6
$ws.Range("L7").Value = 1.2
$ws.Range("N7").Value = 6.2
$ws.Range("O7").Value = 1.15
$ws.Range("P7").Value = 2.82
$ws.Range("Q7").Value = 1.44
$ws.Range("R7").Value = 1.73
$ws.Range("S7").Value = 2.1
$ws.Range("T7").Value = 1.86
$ws.Range("U7").Value = 1.94
$ws.Range("V7").Value = 1.07
$ws.Range("W7").Value = 3.9
$ws.Range("X7").Value = 42
$ws.Range("Y7").Value = 55
$ws.Range("Z7").Value = 140
$ws.Range("AB7").Value = 14.5
$ws.Range("AC7").Value = 19
$ws.Range("AD7").Value = 55
$ws.Range("AF7").Value = 12
$ws.Range("AG7").Value = 14
$ws.Range("AH7").Value = 34
$ws.Range("AJ7").Value = 13
$ws.Range("AK7").Value = 16.5
$ws.Range("AL7").Value = 40
$ws.Range("AN7").Value = 4.6
$ws.Range("I8").Value = 4.8
$ws.Range("V8").Value = 1.26
$ws.Range("X8").Value = 980
$ws.Range("AH8").Value = 980
$ws.Range("AJ8").Value = 980
$ws.Range("AK8").Value = 980
$ws.Range("AL8").Value = 980
$ws.Range("P9").Value = 2.02
$ws.Range("T9").Value = 1.71
$ws.Range("V9").Value = 1.83
$ws.Range("F10").Value = 3.7
$ws.Range("I10").Value = 2.08
$ws.Range("J10").Value = 3.4
$ws.Range("V10").Value = 1.92
$ws.Range("N11").Value = 5.4
$ws.Range("O11").Value = 1.18
$ws.Range("R11").Value = 1.61
$ws.Range("S11").Value = 2.3
$ws.Range("T11").Value = 2.06
$ws.Range("U11").Value = 1.76
$ws.Range("Y11").Value = 60
$ws.Range("Z11").Value = 170
$ws.Range("AB11").Value = 12.5
$ws.Range("AC11").Value = 19
$ws.Range("AD11").Value = 55
$ws.Range("AE11").Value = 280
$ws.Range("AF11").Value = 10.5
$ws.Range("AG11").Value = 14
$ws.Range("AH11").Value = 40
$ws.Range("AI11").Value = 210
$ws.Range("AJ11").Value = 12
$ws.Range("AK11").Value = 17.5
$ws.Range("AL11").Value = 48
$ws.Range("AM11").Value = 210
$ws.Range("AN11").Value = 5
$ws.Range("Q12").Value = 1.9
$ws.Range("R12").Value = 1.16
$ws.Range("S12").Value = 1.9
$ws.Range("F14").Value = 2.22
$ws.Range("J14").Value = 3.4
$ws.Range("L14").Value = 1.24
$ws.Range("N14").Value = 5
$ws.Range("O14").Value = 1.19
$ws.Range("T14").Value = 1.53
$ws.Range("U14").Value = 2.48
$ws.Range("V14").Value = 1.44
$ws.Range("W14").Value = 1.66
$ws.Range("X14").Value = 29
$ws.Range("Y14").Value = 21
$ws.Range("Z14").Value = 30
$ws.Range("AA14").Value = 60
$ws.Range("AB14").Value = 18
$ws.Range("AC14").Value = 10.5
$ws.Range("AD14").Value = 17
$ws.Range("AE14").Value = 32
$ws.Range("AF14").Value = 23
$ws.Range("AG14").Value = 14.5
$ws.Range("AH14").Value = 18.5
$ws.Range("AI14").Value = 42
$ws.Range("AJ14").Value = 38
$ws.Range("AK14").Value = 27
$ws.Range("AL14").Value = 36
$ws.Range("AM14").Value = 70
$ws.Range("AN14").Value = 15.5
$ws.Range("AO14").Value = 24
$ws.Range("G15").Value = 11
$ws.Range("P15").Value = 1.89
$ws.Range("G16").Value = 4.9
$ws.Range("K16").Value = 5
$ws.Range("S16").Value = 2.14
$ws.Range("P18").Value = 1.68
$ws.Range("T18").Value = 1.89
$ws.Range("F19").Value = 2.58
$ws.Range("H19").Value = 3.05
$ws.Range("I19").Value = 3.1
$ws.Range("V19").Value = 1.47
$ws.Range("W19").Value = 1.62
$ws.Range("Z19").Value = 19
$ws.Range("AC19").Value = 7.4
$ws.Range("W21").Value = 3.8
